$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.416.79"
$ws.Range("E2").Value = "  +0.92%  "

# Row 3
$ws.Range("D3").Value = "1.941.85"
$ws.Range("E3").Value = "  -1.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.96"
$ws.Range("E5").Value = "  +0.40%  "

# Row 6
$ws.Range("E6").Value = "  -1.81%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.03"
$ws.Range("E7").Value = "  -6.62%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.364"
$ws.Range("E9").Value = "  -2.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.69"
$ws.Range("E10").Value = "  -0.77%  "

# Row 11
$ws.Range("E11").Value = "  +3.11%  "

# Row 12
$ws.Range("E12").Value = "  +0.62%  "

# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.46"
$ws.Range("E13").Value = "  -2.60%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.817"
$ws.Range("E14").Value = "  -4.79%  "

# Row 15
$ws.Range("D15").Value = "2.227.01"
$ws.Range("E15").Value = "  -1.03%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.49"
$ws.Range("E16").Value = "  -3.83%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.22"
$ws.Range("E17").Value = "  -3.65%  "

# Row 18
$ws.Range("D18").Value = "1.945.02"
$ws.Range("E18").Value = "  -1.10%  "

# Row 19
$ws.Range("D19").Value = "36.270.84"
$ws.Range("E19").Value = "  +0.88%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.62"
$ws.Range("E20").Value = "  -2.05%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0859"
$ws.Range("E21").Value = "  +0.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.08"
$ws.Range("E22").Value = "  -3.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.04"
$ws.Range("E23").Value = "  -2.93%  "

# Row 24
$ws.Range("E24").Value = "  +0.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.42"
$ws.Range("E25").Value = "  -4.85%  "

# Row 26
$ws.Range("E26").Value = "  -0.23%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.19"
$ws.Range("E27").Value = "  -5.96%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.16"
$ws.Range("E28").Value = "  +0.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.36"
$ws.Range("E29").Value = "  -2.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").Value = "  -6.15%  "

# Row 31
$ws.Range("E31").Value = "  -1.62%  "

# Row 32
$ws.Range("E32").Value = "  +0.66%  "

# Row 33
$ws.Range("E33").Value = "  -3.81%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0624"
$ws.Range("E34").Value = "  +0.65%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.23"
$ws.Range("E35").Value = "  -3.94%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.20"
$ws.Range("E36").Value = "  -1.65%  "

# Row 37
$ws.Range("E37").Value = "  +0.09%  "

# Row 38
$ws.Range("E38").Value = "  -2.54%  "

# Row 39
$ws.Range("E39").Value = "  -6.77%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.01"
$ws.Range("E40").Value = "  -2.47%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0984"
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("E42").Value = "  +1.50%  "

# Row 43
$ws.Range("E43").Value = "  -4.33%  "

# Row 44
$ws.Range("E44").Value = "  -1.58%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.98"
$ws.Range("E45").Value = "  +0.13%  "

# Row 46
$ws.Range("D46").Value = "1.345.55"
$ws.Range("E46").Value = "  +0.65%  "

# Row 47
$ws.Range("E47").Value = "  -5.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.47"
$ws.Range("E48").Value = "  -5.23%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.10"
$ws.Range("E49").Value = "  -5.81%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.81"
$ws.Range("E50").Value = "  +0.76%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.21"
$ws.Range("E51").Value = "  +3.30%  "
